$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the first 11 data rows (old rows 2-12), shifting old rows 13-22
# up to become the new rows 2-11.
$ws.Rows("2:12").Delete()

# Append the newly recorded data rows (new rows 12-21) with the
# struggle / walkingToRunning samples added in this commit.
$newData = @(
    @(-2.55821630358695, -8.715493917465219, -6.844864040613141),
    @(-0.8884068131446955, -10.84955549240112, 1.661297619342809),
    @(-3.758049488067625, -8.76593214273451, 3.196124792099005),
    @(-3.226068019866951, -4.671105861663801, 4.809621334075939),
    @(-4.81333899497986, -0.7696201205253579, 7.391847074031832),
    @(-0.6115292310714668, -3.946480035781865, 7.73033595085144),
    @(-4.580467939376843, -2.483635365962976, 9.400022864341736),
    @(-3.186845898628217, 3.314929008483903, 7.02214622497557),
    @(-0.8589091598987735, 5.796535491943349, 4.699775040149696),
    @(-6.590956926345826, -0.2845994234085083, 5.918066263198853)
)

$startRow = 12
for ($i = 0; $i -lt $newData.Count; $i++) {
    $row = $startRow + $i
    $values = $newData[$i]
    $ws.Cells.Item($row, 1).Value = $values[0]
    $ws.Cells.Item($row, 2).Value = $values[1]
    $ws.Cells.Item($row, 3).Value = $values[2]
}
